$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column D
$ws.Range("D2").Value = "canonical SMILES"

# Fill column D ("canonical SMILES") for each microstate row.
# For most rows this equals the existing "canonical isomeric SMILES" value in column C,
# except for micro002 and micro006 where stereo bond markers are stripped.
$ws.Range("D3").Value = "c1ccc(c(c1)NC(=O)c2ccc(o2)Cl)N3CCCCC3"
$ws.Range("D4").Value = "c1ccc(c(c1)N=C(c2ccc(o2)Cl)[O-])N3CCCCC3"
$ws.Range("D5").Value = "c1ccc(c(c1)[NH2+]C(=[OH+])c2ccc(o2)Cl)N3CCCCC3"
$ws.Range("D6").Value = "c1ccc(c(c1)NC(=[OH+])c2ccc(o2)Cl)[NH+]3CCCCC3"
$ws.Range("D7").Value = "c1ccc(c(c1)N=C(c2ccc(o2)Cl)[O-])[NH+]3CCCCC3"
$ws.Range("D8").Value = "c1ccc(c(c1)[N-]C(=[OH+])c2ccc(o2)Cl)N3CCCCC3"
$ws.Range("D9").Value = "c1ccc(c(c1)NC(=[OH+])c2ccc(o2)Cl)N3CCCCC3"
$ws.Range("D10").Value = "c1ccc(c(c1)NC(=O)c2ccc(o2)Cl)[NH+]3CCCCC3"
$ws.Range("D11").Value = "c1ccc(c(c1)[N-]C(=[OH+])c2ccc(o2)Cl)[NH+]3CCCCC3"

# Set width of the new column D to match the target layout
$ws.Columns.Item(4).ColumnWidth = 37
